$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.1190320826869504, 0.306821227259698, 0.7527432677738641, 0.4942365360607697, 0, 1.672833113781282)
    3  = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 10.19245300693656, 0, 15.88780690183548)
    4  = @(0.6606524410359556, 1.655778082260271, 0.7527432677738641, 10.19245300693656, 0, 13.26162679800665)
    5  = @(3.286832544864788, 3286.919754855326, 3.537761648806719, 10.19245300693656, 0, 3303.936802055934)
    6  = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 1, 8.974608811992548)
    7  = @(3.286832544864788, 32771568325.09113, 261.3203778131603, 1133.036916526867, 1, 32771569722.73526)
    8  = @(0.04271373187048222, 0.306821227259698, 0.7527432677738641, 0.4942365360607697, 1, 1.596514762964814)
    9  = @(0.6606524410359556, 0.306821227259698, 0.7527432677738641, 10.19245300693656, 0, 11.91266994300607)
    10 = @(1.455362044514542, 0.306821227259698, 0.1494219747398047, 0.4942365360607697, 0, 2.405841782574814)
    11 = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 10.19245300693656, 1, 18.67282528286833)
    12 = @(1.455362044514542, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 1, 7.143138311642302)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($col = 2; $col -le 7; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 2]
    }
}
